$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 363 (pushing the existing rows 363-489
# down to 364-490, growing the used range from T489 to T490).
$ws.Rows(363).Insert()

$ws.Range("A363").Value2 = 9
$ws.Range("B363").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C363").Value2 = "Metropolitana"
$ws.Range("D363").Value2 = 44809
$ws.Range("E363").Value2 = 13
$ws.Range("F363").Value2 = "Fruta"
$ws.Range("G363").Value2 = 100108
$ws.Range("H363").Value2 = "Tropicales y subtropicales"
$ws.Range("I363").Value2 = 100108002
$ws.Range("J363").Value2 = "Mango"
$ws.Range("K363").Value2 = "Sin especificar"
$ws.Range("L363").Value2 = "Primera"
$ws.Range("M363").Value2 = 710
$ws.Range("N363").Value2 = 8500
$ws.Range("O363").Value2 = 9000
$ws.Range("P363").Value2 = 8754
$ws.Range("Q363").Value2 = "$/bandeja 4 kilos"
$ws.Range("R363").Value2 = "Brasil"
$ws.Range("S363").Value2 = 2188
$ws.Range("T363").Value2 = 4
